# Apply the StructureDefinition-age-gender-group.xlsx update:
#  - Metadata sheet: bump Version/Date, fill in Publisher, replace the
#    duplicated "Contact" row with a "Jurisdiction" row, and drop the
#    now-redundant second "Contact" row entirely (rows shift up by one).
#  - Elements sheet: give the root "Extension" row a real Short/Definition
#    instead of the generic placeholder text.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B3").Value = "6.0.0"
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$wsMeta.Range("B9").Value = "Alvearie Team"
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# The old sheet had row 11 as an exact duplicate of row 10 ("Contact" /
# "No display for ContactDetail"). Delete it so every row below shifts up
# by one, taking the sheet from 21 data rows down to 20.
$wsMeta.Rows(11).Delete()

$wsElements = $wb.Worksheets.Item("Elements")

$wsElements.Range("K2").Value = "Age and Gender Groups"
$wsElements.Range("L2").Value = "Standard code for groupings that combine age and gender"
